$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq "Developer, Administrator, Miss Dina Nasr") {
        $cell.Value = "Miss Dina Nasr, Administrator, Developer"
    }
    elseif ($val -eq "Administrator, Miss Dina Nasr") {
        $cell.Value = "Miss Dina Nasr, Administrator"
    }
}
